# 2018/04/25 edit: expected-profit rate for row 3 (C3) was revised from
# 5.20% to 5.31%, which ripples through the dependent interest/ending-
# balance formulas in D3/E3 and the carried-forward opening balance /
# reinvested-principal formulas in A4/B4. Also record the user's final
# cell selection (moved from F3 to C4) as reflected in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the expected profit rate for row 3.
$ws.Range("C3").Value = 0.0531

# Recalculate so every dependent formula cell carries a fresh cached value.
$wb.Application.Calculate()

# Reflect the active-cell selection recorded in the saved sheet view.
$ws.Range("C4").Select()
